$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 52084
$ws.Range("B2").Value = "Valentina Nunes"
$ws.Range("C2").Value = "Engenharia"
$ws.Range("D2").Value = "Outros"
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 45102
$ws.Range("G2").Value = 7218.02

# Row 3
$ws.Range("A3").Value = 34291
$ws.Range("B3").Value = "Dr. Heitor Melo"
$ws.Range("C3").Value = "Vendas"
$ws.Range("D3").Value = "Doença"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 45093
$ws.Range("G3").Value = 10561.49

# Row 4
$ws.Range("A4").Value = 18282
$ws.Range("B4").Value = "Diego Freitas"
$ws.Range("C4").Value = "Jurídico"
$ws.Range("D4").Value = "Outros"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 45079
$ws.Range("G4").Value = 2620.45

# Row 5
$ws.Range("A5").Value = 83712
$ws.Range("B5").Value = "Miguel Lima"
$ws.Range("C5").Value = "Jurídico"
$ws.Range("D5").Value = "Viagem de negócios"
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 45089
$ws.Range("G5").Value = 3679.58

# Row 6
$ws.Range("A6").Value = 20546
$ws.Range("B6").Value = "Luiza Barros"
$ws.Range("C6").Value = "Vendas"
$ws.Range("D6").Value = "Doença"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 45090
$ws.Range("G6").Value = 5731.52

# Row 7
$ws.Range("A7").Value = 43625
$ws.Range("B7").Value = "Marcela da Mota"
$ws.Range("C7").Value = "Atendimento ao Cliente"
$ws.Range("D7").Value = "Problemas pessoais"
$ws.Range("E7").Value = 8
$ws.Range("F7").Value = 45093
$ws.Range("G7").Value = 9660.97

# Row 8
$ws.Range("A8").Value = 31380
$ws.Range("B8").Value = "Emilly Moraes"
$ws.Range("C8").Value = "Vendas"
$ws.Range("D8").Value = "Problemas pessoais"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 45085
$ws.Range("G8").Value = 8807.79

# Row 9
$ws.Range("A9").Value = 77637
$ws.Range("B9").Value = "Dr. Anthony Sales"
$ws.Range("C9").Value = "Recursos Humanos"
$ws.Range("D9").Value = "Consulta médica"
$ws.Range("E9").Value = 8
$ws.Range("F9").Value = 45105
$ws.Range("G9").Value = 10933.94

# Row 10
$ws.Range("A10").Value = 28604
$ws.Range("B10").Value = "Dr. Thiago da Paz"
$ws.Range("C10").Value = "Vendas"
$ws.Range("D10").Value = "Consulta médica"
$ws.Range("E10").Value = 7
$ws.Range("F10").Value = 45083
$ws.Range("G10").Value = 10800.83

# Row 11
$ws.Range("A11").Value = 4615
$ws.Range("B11").Value = "Yasmin Silva"
$ws.Range("C11").Value = "Vendas"
$ws.Range("D11").Value = "Viagem de negócios"
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = 45087
$ws.Range("G11").Value = 8710.77
